$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 310
$ws1.Range("F3").Value = 70
$ws1.Range("F4").Value = 3759
$ws1.Range("F5").Value = 2268
$ws1.Range("F12").Value = 1395
$ws1.Range("F13").Value = 244
$ws1.Range("F14").Value = 2264
$ws1.Range("F15").Value = 162

# Sheet 4: "全部类型" (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 310
$ws4.Range("F3").Value = 70
$ws4.Range("F4").Value = 3759
$ws4.Range("F5").Value = 2268
$ws4.Range("F15").Value = 1395
$ws4.Range("F16").Value = 244
$ws4.Range("F17").Value = 2264
$ws4.Range("F18").Value = 162
